# Updated symbol list on Wed Jan 11 21:30:20 UTC 2023 with GitHub Actions
# Refresh the crypto price / volume snapshot in columns D (Price) and
# E (Volume(1h)); also fix the BKEXToken / KickToken row ordering swap
# (rows 41-42) that came in with this data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E columns hold numeric-looking text ("279.23", "0.44%", ...) that must
# stay plain text (matches the source inlineStr cells) rather than being
# auto-coerced to numbers/percents by COM. Force the Text number format
# before writing each value so it round-trips as a string, then restore the
# cell's original (default/unstyled) look so no stray formatting is added.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "279.23"
Set-TextValue $ws.Range("E2") "0.44%"

Set-TextValue $ws.Range("E3") "1.00%"

Set-TextValue $ws.Range("D4") "4.833"
Set-TextValue $ws.Range("E4") "-0.76%"

Set-TextValue $ws.Range("D5") "0.06379"
Set-TextValue $ws.Range("E5") "0.27%"

Set-TextValue $ws.Range("D6") "7.033"
Set-TextValue $ws.Range("E6") "0.83%"

Set-TextValue $ws.Range("D7") "1.317"
Set-TextValue $ws.Range("E7") "5.27%"

Set-TextValue $ws.Range("D8") "0.8962"
Set-TextValue $ws.Range("E8") "1.70%"

Set-TextValue $ws.Range("D9") "0.1545"
Set-TextValue $ws.Range("E9") "1.37%"

Set-TextValue $ws.Range("D10") "0.06469"
Set-TextValue $ws.Range("E10") "26.77%"

Set-TextValue $ws.Range("D11") "0.07517"
Set-TextValue $ws.Range("E11") "0.03%"

Set-TextValue $ws.Range("D12") "0.02940"
Set-TextValue $ws.Range("E12") "-0.95%"

Set-TextValue $ws.Range("D13") "0.08998"
Set-TextValue $ws.Range("E13") "-0.09%"

Set-TextValue $ws.Range("D14") "0.001574"
Set-TextValue $ws.Range("E14") "0.64%"

Set-TextValue $ws.Range("D15") "0.0006465"
Set-TextValue $ws.Range("E15") "0.62%"

Set-TextValue $ws.Range("D16") "0.006088"
Set-TextValue $ws.Range("E16") "2.37%"

Set-TextValue $ws.Range("E17") "0.71%"

Set-TextValue $ws.Range("D18") "3.327"
Set-TextValue $ws.Range("E18") "0.17%"

Set-TextValue $ws.Range("D19") "2.232"
Set-TextValue $ws.Range("E19") "-1.77%"

Set-TextValue $ws.Range("D21") "0.1352"
Set-TextValue $ws.Range("E21") "1.09%"

Set-TextValue $ws.Range("D22") "3.900"
Set-TextValue $ws.Range("E22") "0.00%"

Set-TextValue $ws.Range("D23") "0.04395"
Set-TextValue $ws.Range("E23") "-0.57%"

Set-TextValue $ws.Range("D24") "0.1503"
Set-TextValue $ws.Range("E24") "8.93%"

Set-TextValue $ws.Range("D25") "0.001175"
Set-TextValue $ws.Range("E25") "0.17%"

Set-TextValue $ws.Range("D26") "0.004280"
Set-TextValue $ws.Range("E26") "10.34%"

Set-TextValue $ws.Range("D28") "0.0001180"
Set-TextValue $ws.Range("E28") "-1.72%"

Set-TextValue $ws.Range("D29") "0.0001653"
Set-TextValue $ws.Range("E29") "-14.59%"

Set-TextValue $ws.Range("D40") "0.04064"
Set-TextValue $ws.Range("E40") "-2.00%"

# Row 41 and 42 swap places: BKEXToken -> KickToken, KickToken -> BKEXToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006672"
Set-TextValue $ws.Range("E41") "-2.32%"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1400"
Set-TextValue $ws.Range("E42") "18.53%"

Set-TextValue $ws.Range("D43") "0.002079"
Set-TextValue $ws.Range("E43") "2.92%"

Set-TextValue $ws.Range("D44") "0.01103"
Set-TextValue $ws.Range("E44") "-1.47%"

Set-TextValue $ws.Range("D45") "0.00005543"
Set-TextValue $ws.Range("E45") "6.90%"

Set-TextValue $ws.Range("D46") "1.561"
Set-TextValue $ws.Range("E46") "5.00%"

Set-TextValue $ws.Range("D47") "0.01849"
Set-TextValue $ws.Range("E47") "-8.68%"
